$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: force text storage for numeric-looking strings, then restore default style

# Row 47 and 48: USDe and Arweave swap positions with updated price/volume values
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.14%  "

$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.02%  "

# Update Price (D) and Volume(1h) (E) columns for all other changed rows
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.114.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.733.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.735.34"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.53%  "
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("E9").Value = "  +2.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.166"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.32"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.463"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000246"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.352.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.723.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.996.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.15%  "
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "498.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.728"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000142"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.72"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.866.02"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.109"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.669.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.134"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.325"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "436.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "49.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.87"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "142.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0352"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.749.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.89%  "
